# Storymap.xlsx update
# The "User" swimlane gains a "Registrieren" (register) flow next to the existing
# "Anmelden" (login) flow; the old inline "...ändern" (change data) steps are
# replaced by a dedicated "Registrieren" column whose sub-steps are
# "Namen/Passwort/Postleitzahl/E-Mail/Adresse eingeben" (enter name/password/
# zip/e-mail/address). Everything to the right of that column shifts from
# column G onward to the next free column (E/I/M instead of G/K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Cells that keep their position/style but receive new text
$ws.Range("C3").Value  = "Registrieren"
$ws.Range("C5").Value  = "Registrieren"
$ws.Range("E5").Value  = "Konto löschen"
$ws.Range("G5").Value  = "Daten ändern"
$ws.Range("I5").Value  = "Artikel in Warenkorb legen"
$ws.Range("K5").Value  = "Bestellung abgeben"
$ws.Range("C7").Value  = "Namen eingeben"
$ws.Range("E7").Value  = "Konto löschen"
$ws.Range("I7").Value  = "Artikel nach Suchbegriff suchen"
$ws.Range("K7").Value  = "Waren in Warenkorb bestellen"

# 2) Brand-new cells: set the text, then copy the cell format (fill/border/
#    alignment) from a same-row sibling so the new cell matches its row style
$ws.Range("E3").Value  = "Konto Verwaltung"
$ws.Range("I3").Value  = "Artikel bestellen"
$ws.Range("M5").Value  = "Bestellung stornieren"
$ws.Range("M7").Value  = "Bestellung auswählen"
$ws.Range("C9").Value  = "Passwort eingeben"
$ws.Range("I9").Value  = "Artikelmenge auswählen"
$ws.Range("M9").Value  = "Bestellung stornieren"
$ws.Range("C11").Value = "Postleitzahl eingeben"
$ws.Range("I11").Value = "Artikel in den Warenkorb legen"
$ws.Range("C13").Value = "E-Mail eingeben"
$ws.Range("C15").Value = "Adresse eingeben"

$ws.Range("A3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Copy() | Out-Null
$ws.Range("M5").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Copy() | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null
$ws.Range("AA9").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("AA9").Copy() | Out-Null
$ws.Range("I9").PasteSpecial(-4122) | Out-Null
$ws.Range("AA9").Copy() | Out-Null
$ws.Range("M9").PasteSpecial(-4122) | Out-Null
$ws.Range("AA11").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("AA11").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null
$ws.Range("W13").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("W15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Cells that no longer exist in the updated layout
$ws.Range("G3").Clear() | Out-Null
$ws.Range("G7").Clear() | Out-Null
$ws.Range("E9").Clear() | Out-Null
$ws.Range("G9").Clear() | Out-Null
$ws.Range("K9").Clear() | Out-Null
$ws.Range("E11").Clear() | Out-Null
$ws.Range("G11").Clear() | Out-Null
$ws.Range("E13").Clear() | Out-Null
$ws.Range("E15").Clear() | Out-Null

# 4) Match the saved selection/active cell
$ws.Range("C3").Select() | Out-Null
